$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.944.72'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.218.49'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '291.26'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").Value = '86.70'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Value = '30.41'
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").Value = '50.27'
$ws.Range("E11").Value = '  +5.83%  '
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Value = '6.43'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").Value = '2.561.12'
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '13.77'
$ws.Range("E16").Value = '  -2.18%  '
$ws.Range("D17").Value = '2.229.34'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").Value = '0.731'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = '39.875.84'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").Value = '11.08'
$ws.Range("E21").Value = '  -3.83%  '
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").Value = '65.67'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '237.97'
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").Value = '23.02'
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  -7.59%  '
$ws.Range("D31").Value = '156.36'
$ws.Range("E31").Value = '  +2.82%  '
$ws.Range("D32").Value = '31.83'
$ws.Range("E32").Value = '  -2.97%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  +6.06%  '
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").Value = '0.0990'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '1.73'
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("D41").Value = '15.22'
$ws.Range("E41").Value = '  -5.00%  '
$ws.Range("D42").Value = '2.095.09'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("D44").Value = '0.0272'
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("D45").Value = '17.95'
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("E46").Value = '  -2.35%  '
$ws.Range("E47").Value = '  -8.12%  '
$ws.Range("E48").Value = '  +3.25%  '
$ws.Range("D49").Value = '2.432.40'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +0.47%  '
